# Fix Training Data Issue (#48)
# * Fix data issue for model training
# * Data was taken from 1 day off due to way NBA stats were shown
#
# The "Date" column held a malformed label ("2-16-2012-13", i.e. the game
# day "2-16" concatenated with the season "2012-13") instead of the actual
# game date. Replace it with the correct ISO-style date string
# "2013-02-16" for every data row, keeping the value as plain text (the
# leading apostrophe stops Excel from reinterpreting the date-looking
# string as a date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Locate the "Date" header column (BF in the known layout) dynamically.
$dateCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value2 -eq "Date") {
        $dateCol = $c
        break
    }
}

if ($dateCol -eq 0) {
    $dateCol = 58  # fallback: column BF
}

$oldValue = "2-16-2012-13"
$newValue = "2013-02-16"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $dateCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = "'" + $newValue
    }
}
